$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.150.76"
$ws.Range("E2").Value = "  -3.30%  "
$ws.Range("D3").Value = "1.713.68"
$ws.Range("E3").Value = "  -3.65%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'309.12"
$ws.Range("E5").Value = "  -6.01%  "
$ws.Range("D6").Value = "'1.000"
$ws.Range("E6").Value = "  +0.10%  "
$ws.Range("D7").Value = "'0.4790"
$ws.Range("E7").Value = "  +6.31%  "
$ws.Range("D8").Value = "'0.3452"
$ws.Range("E8").Value = "  -3.14%  "
$ws.Range("D9").Value = "'42.14"
$ws.Range("E9").Value = "  +0.24%  "
$ws.Range("D10").Value = "'0.07281"
$ws.Range("E10").Value = "  -2.30%  "
$ws.Range("D11").Value = "'1.045"
$ws.Range("E11").Value = "  -5.70%  "
$ws.Range("D12").Value = "'1.000"
$ws.Range("E12").Value = "  +0.06%  "
$ws.Range("D13").Value = "'19.86"
$ws.Range("E13").Value = "  -5.61%  "
$ws.Range("D14").Value = "'5.867"
$ws.Range("E14").Value = "  -3.17%  "
$ws.Range("D15").Value = "1.712.61"
$ws.Range("E15").Value = "  -3.51%  "
$ws.Range("D16").Value = "'6.862"
$ws.Range("E16").Value = "  -5.63%  "
$ws.Range("D17").Value = "'88.94"
$ws.Range("E17").Value = "  -5.47%  "
$ws.Range("E18").Value = "  -2.22%  "
$ws.Range("D19").Value = "'0.06367"
$ws.Range("E19").Value = "  -1.24%  "
$ws.Range("D20").Value = "'1.000"
$ws.Range("E20").Value = "  +0.11%  "
$ws.Range("D21").Value = "'16.51"
$ws.Range("E21").Value = "  -3.69%  "
$ws.Range("D22").Value = "'5.622"
$ws.Range("E22").Value = "  -3.11%  "
$ws.Range("D23").Value = "27.183.32"
$ws.Range("E23").Value = "  -3.19%  "
$ws.Range("D24").Value = "'10.84"
$ws.Range("E24").Value = "  -4.32%  "
$ws.Range("D25").Value = "'2.088"
$ws.Range("E25").Value = "  -1.84%  "
$ws.Range("D26").Value = "'151.92"
$ws.Range("E26").Value = "  -6.12%  "
$ws.Range("D27").Value = "'19.69"
$ws.Range("E27").Value = "  -3.54%  "
$ws.Range("D28").Value = "1.908.75"
$ws.Range("E28").Value = "  -3.57%  "
$ws.Range("D29").Value = "'2.092"
$ws.Range("E29").Value = "  -3.33%  "
$ws.Range("D30").Value = "'120.06"
$ws.Range("E30").Value = "  -3.94%  "
$ws.Range("D31").Value = "'1.018"
$ws.Range("E31").Value = "  -8.27%  "
$ws.Range("D32").Value = "'0.09282"
$ws.Range("E32").Value = "  +0.66%  "
$ws.Range("D33").Value = "'3.585"
$ws.Range("E33").Value = "  -3.05%  "
$ws.Range("D34").Value = "'5.315"
$ws.Range("E34").Value = "  -7.23%  "
$ws.Range("D35").Value = "'0.02203"
$ws.Range("E35").Value = "  -4.03%  "
$ws.Range("D36").Value = "'0.05908"
$ws.Range("E36").Value = "  -4.92%  "
$ws.Range("D37").Value = "'11.07"
$ws.Range("E37").Value = "  -6.97%  "
$ws.Range("D38").Value = "'0.2006"
$ws.Range("E38").Value = "  -5.14%  "
$ws.Range("B39").Value = "WEMIXTOKEN"
$ws.Range("C39").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D39").Value = "'1.418"
$ws.Range("E39").Value = "  +1.41%  "
$ws.Range("B40").Value = "InternetComputer(DFINITY)"
$ws.Range("C40").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D40").Value = "'4.748"
$ws.Range("E40").Value = "  -5.04%  "
$ws.Range("E41").Value = "  +0.12%  "
$ws.Range("D42").Value = "'0.5936"
$ws.Range("E42").Value = "  -6.21%  "
$ws.Range("E43").Value = "  -6.38%  "
$ws.Range("D44").Value = "'7.479"
$ws.Range("E44").Value = "  -5.93%  "
$ws.Range("D45").Value = "'12.74"
$ws.Range("E45").Value = "  -4.10%  "
$ws.Range("D46").Value = "'3.571"
$ws.Range("E46").Value = "  -4.95%  "
$ws.Range("D47").Value = "'0.5623"
$ws.Range("E47").Value = "  -4.86%  "
$ws.Range("D48").Value = "'118.76"
$ws.Range("E48").Value = "  -3.27%  "
$ws.Range("D49").Value = "'1.841"
$ws.Range("D50").Value = "'0.06642"
$ws.Range("E50").Value = "  -3.74%  "
$ws.Range("D51").Value = "'1.084"
$ws.Range("E51").Value = "  -5.12%  "
